$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44677; J = 20;  K = 5500; L = 5500; M = 5500; P = 5500 }
    3  = @{ D = 44719; J = 80;  K = 3600; L = 3600; M = 3600; P = 3600 }
    4  = @{ D = 44203; J = 30;  K = 2000; L = 2000; M = 2000; P = 2000 }
    5  = @{ D = 44669; J = 60;  K = 6250; L = 6250; M = 6250; P = 6250 }
    6  = @{ D = 44497; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
    7  = @{ D = 44679; J = 30;  K = 5500; L = 5500; M = 5500; P = 5500 }
    8  = @{ D = 44447; J = 75;  K = 2200; L = 2200; M = 2200; P = 2200 }
    9  = @{ D = 44453; J = 20;  K = 2300; L = 2300; M = 2300; P = 2300 }
    10 = @{ D = 44487; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
    11 = @{ D = 44496; J = 40;  K = 2200; L = 2200; M = 2200; P = 2200 }
    12 = @{ D = 44484; J = 40;  K = 2200; L = 2200; M = 2200; P = 2200 }
    13 = @{ D = 44685; J = 60;  K = 5000; L = 6000; M = 5333; P = 5333 }
    14 = @{ D = 44720; J = 100; K = 3600; L = 3600; M = 3600; P = 3600 }
    15 = @{ D = 44452; J = 120; K = 2300; L = 2300; M = 2300; P = 2300 }
    16 = @{ D = 44706; J = 90;  K = 4700; L = 4700; M = 4700; P = 4700 }
    17 = @{ D = 44476; J = 30;  K = 2200; L = 2200; M = 2200; P = 2200 }
    18 = @{ D = 44474; J = 20;  K = 1600; L = 1600; M = 1600; P = 1600 }
    20 = @{ D = 44707; J = 100; K = 4700; L = 4700; M = 4700; P = 4700 }
    21 = @{ D = 44483; J = 50;  K = 2200; L = 2200; M = 2200; P = 2200 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
